# Adds this week's price report (week of 2021-11-22, serial 44522) for
# "Agrícola del Norte S.A. de Arica - Pepino ensalada".
# The existing historical rows (186:207) are pushed down by two rows
# (188:209) and the freed-up rows 186:187 are filled with the new
# Primera/Segunda quality records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data block down by inserting two fresh rows at the
# top of it (mirrors Excel's "Insert Copied/Blank Cells" row behaviour -
# it also carries the date-column style down onto the new rows).
$ws.Rows("186:187").Insert()

# Row 186 - Primera
$ws.Cells.Item(186, 1).Value = 1
$ws.Cells.Item(186, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(186, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(186, 4).Value = 44522
$ws.Cells.Item(186, 5).Value = 15
$ws.Cells.Item(186, 6).Value = 100112043
$ws.Cells.Item(186, 7).Value = "Pepino ensalada"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 120
$ws.Cells.Item(186, 11).Value = 4000
$ws.Cells.Item(186, 12).Value = 4500
$ws.Cells.Item(186, 13).Value = 4250
$ws.Cells.Item(186, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(186, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(186, 16).Value = 61
$ws.Cells.Item(186, 17).Value = 70
$ws.Cells.Item(186, 18).Value = "Hortaliza"

# Row 187 - Segunda
$ws.Cells.Item(187, 1).Value = 1
$ws.Cells.Item(187, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(187, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(187, 4).Value = 44522
$ws.Cells.Item(187, 5).Value = 15
$ws.Cells.Item(187, 6).Value = 100112043
$ws.Cells.Item(187, 7).Value = "Pepino ensalada"
$ws.Cells.Item(187, 8).Value = "Sin especificar"
$ws.Cells.Item(187, 9).Value = "Segunda"
$ws.Cells.Item(187, 10).Value = 140
$ws.Cells.Item(187, 11).Value = 3000
$ws.Cells.Item(187, 12).Value = 3500
$ws.Cells.Item(187, 13).Value = 3250
$ws.Cells.Item(187, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(187, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(187, 16).Value = 32
$ws.Cells.Item(187, 17).Value = 100
$ws.Cells.Item(187, 18).Value = "Hortaliza"
